# Corrections in Correlation analysis
# Update the "Bugs" counts (column M) and re-point the ratio formulas in
# column B at them (one denominator also changed for row 4), then move the
# active selection from B6 to B4 to match where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- corrected bug counts (M column) ---
$ws.Range("M2").Value2 = 117
$ws.Range("M3").Value2 = 50
$ws.Range("M4").Value2 = 17
# M5 and M6 are unchanged (49 and 111 respectively)

# --- formulas now reference the M column instead of hard-coded numerators ---
$ws.Range("B2").Formula = "=M2/296"
$ws.Range("B3").Formula = "=M3/66"
$ws.Range("B4").Formula = "=M4/127"
$ws.Range("B5").Formula = "=M5/297"
$ws.Range("B6").Formula = "=M6/38"

# --- move the active selection to B4 ---
$ws.Range("B4").Select()
